# Auto-update draw results: append the 2025-11-03 Pick 3 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every column as plain text (dates, phase codes, and
# results are all literal strings, not real dates/numbers). Force the new
# row to Text format first so Excel doesn't "helpfully" convert
# "2025-11-03" into a date serial or "251103" into a number.
$row = $ws.Range("A48:E48")
$row.NumberFormat = "@"

$ws.Range("A48").Value = "2025-11-03"
$ws.Range("B48").Value = "Pick 3"
$ws.Range("C48").Value = "251103"
$ws.Range("D48").Value = "1-4-3"
$ws.Range("E48").Value = "2025-11-03T21:39:02.052+04:00"
